# The deck's cached "date" placeholder (the literal text cached inside the
# datetimeFigureOut / datetime1 <a:fld> elements) needs to move forward one
# week: 2026-01-16 -> 2026-01-23. That placeholder lives on the slide
# master, on every one of its slide layouts, and on the notes master.

$p = $ppt.ActivePresentation

$oldDate = "2026-01-16"
$newDate = "2026-01-23"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    $count = $shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $isDatePh = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
            if ($isDatePh -eq $true) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every slide layout (CustomLayouts) hanging off the slide master.
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j)
}
